$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.576.30"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "3.475.05"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.37"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("E7").Value = "  -3.00%  "
$ws.Range("D8").Value = "3.466.92"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.77"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.91%  "
$ws.Range("E12").Value = "  -5.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.42"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "4.032.76"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "614.05"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -11.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.29"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.19%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.478.19"
$ws.Range("E18").Value = "  -2.92%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "68.614.47"
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.20"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.871"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.70"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.57"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -4.29%  "
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.72"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.38"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("E33").Value = "  -3.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.80"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "576.54"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.66"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.47"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -10.10%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.69"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0434"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.136"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("D43").Value = "3.400.80"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.56"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.22%  "
$ws.Range("D46").Value = "0.0₃0688"
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.79"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("E49").Value = "  -4.18%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +11.72%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.32"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.38%  "
